$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (41 cell updates) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4739.5557
$ws.Cells.Item(32, 9).Value = 2516.8
$ws.Cells.Item(32, 10).Value = 5594.4614
$ws.Cells.Item(32, 11).Value = 2516.8
$ws.Cells.Item(32, 12).Value = 5594.4614
$ws.Cells.Item(32, 13).Value = -2190.8
$ws.Cells.Item(32, 14).Value = -6246.4614
$ws.Cells.Item(88, 8).Value = 20009570
$ws.Cells.Item(88, 9).Value = 66667570
$ws.Cells.Item(88, 10).Value = 13284.857
$ws.Cells.Item(88, 11).Value = 66667570
$ws.Cells.Item(88, 12).Value = 13284.857
$ws.Cells.Item(88, 13).Value = -66667164
$ws.Cells.Item(88, 14).Value = -14096.857
$ws.Cells.Item(91, 8).Value = 20009570
$ws.Cells.Item(91, 9).Value = 66667570
$ws.Cells.Item(91, 10).Value = 13284.857
$ws.Cells.Item(91, 11).Value = 66667570
$ws.Cells.Item(91, 12).Value = 13284.857
$ws.Cells.Item(91, 13).Value = -66666166
$ws.Cells.Item(91, 14).Value = -16092.857
$ws.Cells.Item(98, 8).Value = 3586.25
$ws.Cells.Item(98, 9).Value = 3166.5
$ws.Cells.Item(98, 11).Value = 3166.5
$ws.Cells.Item(98, 13).Value = -1668.5
$ws.Cells.Item(113, 8).Value = 3928.1428
$ws.Cells.Item(113, 10).Value = 3999.4
$ws.Cells.Item(113, 12).Value = 3999.4
$ws.Cells.Item(113, 14).Value = -10507.4
$ws.Cells.Item(122, 8).Value = 3586.25
$ws.Cells.Item(122, 9).Value = 3166.5
$ws.Cells.Item(122, 11).Value = 9499.5
$ws.Cells.Item(122, 13).Value = -7049.5
$ws.Cells.Item(125, 8).Value = 2847.8572
$ws.Cells.Item(125, 9).Value = 858.75
$ws.Cells.Item(125, 11).Value = 7728.75
$ws.Cells.Item(125, 13).Value = -5268.75
$ws.Cells.Item(127, 8).Value = 2590.1191
$ws.Cells.Item(127, 9).Value = 911.2727
$ws.Cells.Item(127, 11).Value = 2733.8181
$ws.Cells.Item(127, 13).Value = 2226.1819

# ---- Sheet: ARM (19 cell updates) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 796.8333
$ws.Cells.Item(4, 9).Value = 194.75
$ws.Cells.Item(4, 10).Value = 2001
$ws.Cells.Item(4, 11).Value = 194.75
$ws.Cells.Item(4, 12).Value = 2001
$ws.Cells.Item(4, 13).Value = -78.75
$ws.Cells.Item(4, 14).Value = -2233
$ws.Cells.Item(74, 8).Value = 1357.7273
$ws.Cells.Item(74, 9).Value = 1048.5555
$ws.Cells.Item(74, 11).Value = 1048.5555
$ws.Cells.Item(74, 13).Value = -174.5554999999999
$ws.Cells.Item(77, 8).Value = 1357.7273
$ws.Cells.Item(77, 9).Value = 1048.5555
$ws.Cells.Item(77, 11).Value = 5242.7775
$ws.Cells.Item(77, 13).Value = -874.7775000000001
$ws.Cells.Item(109, 8).Value = 89249.5
$ws.Cells.Item(109, 10).Value = 89249.5
$ws.Cells.Item(109, 12).Value = 89249.5
$ws.Cells.Item(109, 14).Value = -92023.5

# ---- Sheet: BSM (12 cell updates) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3054.48
$ws.Cells.Item(86, 9).Value = 1960.4375
$ws.Cells.Item(86, 11).Value = 1960.4375
$ws.Cells.Item(86, 13).Value = -837.4375
$ws.Cells.Item(89, 8).Value = 3054.48
$ws.Cells.Item(89, 9).Value = 1960.4375
$ws.Cells.Item(89, 11).Value = 9802.1875
$ws.Cells.Item(89, 13).Value = -4186.1875
$ws.Cells.Item(134, 8).Value = 3397.1333
$ws.Cells.Item(134, 9).Value = 2549.2727
$ws.Cells.Item(134, 11).Value = 7647.8181
$ws.Cells.Item(134, 13).Value = -5112.8181

# ---- Sheet: CRP (72 cell updates) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 3054.8276
$ws.Cells.Item(7, 9).Value = 2989.8948
$ws.Cells.Item(7, 10).Value = 3178.2
$ws.Cells.Item(7, 11).Value = 2989.8948
$ws.Cells.Item(7, 12).Value = 3178.2
$ws.Cells.Item(7, 13).Value = -2876.8948
$ws.Cells.Item(7, 14).Value = -3404.2
$ws.Cells.Item(16, 8).Value = 1418.2858
$ws.Cells.Item(16, 9).Value = 965
$ws.Cells.Item(16, 10).Value = 2022.6666
$ws.Cells.Item(16, 11).Value = 965
$ws.Cells.Item(16, 12).Value = 2022.6666
$ws.Cells.Item(16, 13).Value = -678
$ws.Cells.Item(16, 14).Value = -2596.6666
$ws.Cells.Item(31, 8).Value = 13700726
$ws.Cells.Item(31, 10).Value = 6619
$ws.Cells.Item(31, 12).Value = 6619
$ws.Cells.Item(31, 14).Value = -7209
$ws.Cells.Item(34, 8).Value = 13700726
$ws.Cells.Item(34, 10).Value = 6619
$ws.Cells.Item(34, 12).Value = 6619
$ws.Cells.Item(34, 14).Value = -7023
$ws.Cells.Item(58, 8).Value = 3320.111
$ws.Cells.Item(58, 9).Value = 2766.7144
$ws.Cells.Item(58, 10).Value = 5257
$ws.Cells.Item(58, 11).Value = 2766.7144
$ws.Cells.Item(58, 12).Value = 5257
$ws.Cells.Item(58, 13).Value = -2563.7144
$ws.Cells.Item(58, 14).Value = -5663
$ws.Cells.Item(104, 8).Value = 49999.668
$ws.Cells.Item(104, 9).Value = 59999
$ws.Cells.Item(104, 10).Value = 45000
$ws.Cells.Item(104, 11).Value = 59999
$ws.Cells.Item(104, 12).Value = 45000
$ws.Cells.Item(104, 13).Value = -57378
$ws.Cells.Item(104, 14).Value = -50242
$ws.Cells.Item(105, 8).Value = 1624197.4
$ws.Cells.Item(105, 9).Value = 2066705.9
$ws.Cells.Item(105, 11).Value = 2066705.9
$ws.Cells.Item(105, 13).Value = -2064958.9
$ws.Cells.Item(113, 8).Value = 1418.2858
$ws.Cells.Item(113, 9).Value = 965
$ws.Cells.Item(113, 10).Value = 2022.6666
$ws.Cells.Item(113, 11).Value = 965
$ws.Cells.Item(113, 12).Value = 2022.6666
$ws.Cells.Item(113, 13).Value = 1205
$ws.Cells.Item(113, 14).Value = -6362.6666
$ws.Cells.Item(122, 8).Value = 2998.6667
$ws.Cells.Item(122, 9).Value = 2098.4
$ws.Cells.Item(122, 10).Value = 7500
$ws.Cells.Item(122, 11).Value = 6295.200000000001
$ws.Cells.Item(122, 12).Value = 22500
$ws.Cells.Item(122, 13).Value = -3845.200000000001
$ws.Cells.Item(122, 14).Value = -27400
$ws.Cells.Item(132, 8).Value = 88902060
$ws.Cells.Item(132, 9).Value = 95245064
$ws.Cells.Item(132, 10).Value = 100000
$ws.Cells.Item(132, 11).Value = 285735192
$ws.Cells.Item(132, 12).Value = 300000
$ws.Cells.Item(132, 13).Value = -285732662
$ws.Cells.Item(132, 14).Value = -305060
$ws.Cells.Item(136, 8).Value = 3320.111
$ws.Cells.Item(136, 9).Value = 2766.7144
$ws.Cells.Item(136, 10).Value = 5257
$ws.Cells.Item(136, 11).Value = 8300.143199999999
$ws.Cells.Item(136, 12).Value = 15771
$ws.Cells.Item(136, 13).Value = -5750.143199999999
$ws.Cells.Item(136, 14).Value = -20871
$ws.Cells.Item(141, 8).Value = 114074.73
$ws.Cells.Item(141, 10).Value = 114936.17
$ws.Cells.Item(141, 12).Value = 114936.17
$ws.Cells.Item(141, 14).Value = -125296.17

# ---- Sheet: CUL (4 cell updates) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 6161
$ws.Cells.Item(52, 10).Value = 6161
$ws.Cells.Item(52, 12).Value = 18483
$ws.Cells.Item(52, 14).Value = -19015

# ---- Sheet: GSM (34 cell updates) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 925983.4399999999
$ws.Cells.Item(80, 10).Value = 2963.75
$ws.Cells.Item(80, 12).Value = 2963.75
$ws.Cells.Item(80, 14).Value = -4959.75
$ws.Cells.Item(83, 8).Value = 925983.4399999999
$ws.Cells.Item(83, 10).Value = 2963.75
$ws.Cells.Item(83, 12).Value = 14818.75
$ws.Cells.Item(83, 14).Value = -24802.75
$ws.Cells.Item(97, 8).Value = 685.64703
$ws.Cells.Item(97, 9).Value = 647.9167
$ws.Cells.Item(97, 10).Value = 776.2
$ws.Cells.Item(97, 11).Value = 647.9167
$ws.Cells.Item(97, 12).Value = 776.2
$ws.Cells.Item(97, 13).Value = -151.9167
$ws.Cells.Item(97, 14).Value = -1768.2
$ws.Cells.Item(102, 8).Value = 7861.1665
$ws.Cells.Item(102, 9).Value = 7761.1113
$ws.Cells.Item(102, 11).Value = 7761.1113
$ws.Cells.Item(102, 13).Value = -6139.1113
$ws.Cells.Item(122, 8).Value = 594631.7
$ws.Cells.Item(122, 9).Value = 1432035.6
$ws.Cells.Item(122, 11).Value = 4296106.800000001
$ws.Cells.Item(122, 13).Value = -4293656.800000001
$ws.Cells.Item(123, 8).Value = 33991.668
$ws.Cells.Item(123, 10).Value = 33991.668
$ws.Cells.Item(123, 12).Value = 33991.668
$ws.Cells.Item(123, 14).Value = -38891.668
$ws.Cells.Item(126, 8).Value = 4872.5625
$ws.Cells.Item(126, 9).Value = 2104
$ws.Cells.Item(126, 10).Value = 6533.7
$ws.Cells.Item(126, 11).Value = 6312
$ws.Cells.Item(126, 12).Value = 19601.1
$ws.Cells.Item(126, 13).Value = -3842
$ws.Cells.Item(126, 14).Value = -24541.1

# ---- Sheet: LTW (40 cell updates) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1046.1428
$ws.Cells.Item(22, 9).Value = 899.625
$ws.Cells.Item(22, 10).Value = 1241.5
$ws.Cells.Item(22, 11).Value = 899.625
$ws.Cells.Item(22, 12).Value = 1241.5
$ws.Cells.Item(22, 13).Value = -604.625
$ws.Cells.Item(22, 14).Value = -1831.5
$ws.Cells.Item(27, 8).Value = 1046.1428
$ws.Cells.Item(27, 9).Value = 899.625
$ws.Cells.Item(27, 10).Value = 1241.5
$ws.Cells.Item(27, 11).Value = 899.625
$ws.Cells.Item(27, 12).Value = 1241.5
$ws.Cells.Item(27, 13).Value = -792.625
$ws.Cells.Item(27, 14).Value = -1455.5
$ws.Cells.Item(40, 8).Value = 9619484
$ws.Cells.Item(40, 9).Value = 13891810
$ws.Cells.Item(40, 10).Value = 6750
$ws.Cells.Item(40, 11).Value = 13891810
$ws.Cells.Item(40, 12).Value = 6750
$ws.Cells.Item(40, 13).Value = -13891674
$ws.Cells.Item(40, 14).Value = -7022
$ws.Cells.Item(55, 8).Value = 337.4737
$ws.Cells.Item(55, 9).Value = 72.333336
$ws.Cells.Item(55, 11).Value = 72.333336
$ws.Cells.Item(55, 13).Value = 100.666664
$ws.Cells.Item(93, 8).Value = 1587.7222
$ws.Cells.Item(93, 9).Value = 1428.9286
$ws.Cells.Item(93, 11).Value = 1428.9286
$ws.Cells.Item(93, 13).Value = -180.9286
$ws.Cells.Item(100, 8).Value = 1995.4286
$ws.Cells.Item(100, 9).Value = 1992.25
$ws.Cells.Item(100, 10).Value = 1999.6666
$ws.Cells.Item(100, 11).Value = 1992.25
$ws.Cells.Item(100, 12).Value = 1999.6666
$ws.Cells.Item(100, 13).Value = -1451.25
$ws.Cells.Item(100, 14).Value = -3081.6666
$ws.Cells.Item(141, 8).Value = 98000
$ws.Cells.Item(141, 10).Value = 98000
$ws.Cells.Item(141, 12).Value = 98000
$ws.Cells.Item(141, 14).Value = -108360

# ---- Sheet: WVR (23 cell updates) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(7, 8).Value = 1004
$ws.Cells.Item(7, 9).Value = 1004
$ws.Cells.Item(7, 11).Value = 1004
$ws.Cells.Item(7, 13).Value = -891
$ws.Cells.Item(122, 8).Value = 3441.0386
$ws.Cells.Item(122, 9).Value = 2073.35
$ws.Cells.Item(122, 11).Value = 6220.049999999999
$ws.Cells.Item(122, 13).Value = -3770.049999999999
$ws.Cells.Item(132, 8).Value = 2871.4783
$ws.Cells.Item(132, 9).Value = 2638.3635
$ws.Cells.Item(132, 10).Value = 8000
$ws.Cells.Item(132, 11).Value = 7915.0905
$ws.Cells.Item(132, 12).Value = 24000
$ws.Cells.Item(132, 13).Value = -5385.0905
$ws.Cells.Item(132, 14).Value = -29060
$ws.Cells.Item(138, 8).Value = 126699.8
$ws.Cells.Item(138, 10).Value = 135950
$ws.Cells.Item(138, 12).Value = 135950
$ws.Cells.Item(138, 14).Value = -146230
$ws.Cells.Item(140, 8).Value = 85063.75
$ws.Cells.Item(140, 10).Value = 85063.75
$ws.Cells.Item(140, 12).Value = 85063.75
$ws.Cells.Item(140, 14).Value = -95423.75
